# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> bound to the slide master (drives every slide)
#     originally held the "Integral" color scheme
#   ppt/theme/theme2.xml -> bound to the notes master
#     originally held the stock "Office Theme" color scheme
#
# The authored edit swaps the two themes' palettes (the font scheme and
# format/fill/line/effect scheme are byte-identical between the two themes,
# so only the 12 theme colors actually move). Re-create that by pushing the
# "Office Theme" color scheme values onto the live ThemeColorScheme that
# backs theme1.xml (the slide master's theme), mirroring what Office does
# when a new color scheme is applied to the deck.

function RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order is document order of <a:clrScheme>: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$tcs.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink 954F72
